$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: "หัวcard" column (B) changes from the "unlimited internet" promo
# to the new "โปรฮิต สุดคุ้ม" promo name.
$ws.Range("B2").Value = "โปรฮิต สุดคุ้ม"

# Row 3: drop the card-name value entirely (cell becomes empty/absent),
# and swap the promo-type / bonus columns.
$ws.Range("B3").ClearContents()
$ws.Range("C3").Value = "เน็ตไม่จำกัด"
$ws.Range("G3").Value = "300 บาท / 30วัน"
$ws.Range("H3").Value = "AIS SUPER WiFi ไม่จำกัด"

# Row 4: rename the card and update the promo details to match the
# new 300-baht / 30-day package.
$ws.Range("B4").Value = "โปรฮิต สุดคุ้ม"
$ws.Range("C4").Value = "เน็ตไม่จำกัด"
$ws.Range("E4").Value = "300 บาท"
$ws.Range("F4").Value = "30 วัน"
$ws.Range("G4").Value = "300 บาท / 30วัน"

# Move the active selection the way the author's session left it.
$ws.Range("F13").Select() | Out-Null
